# Natmi following Dr Hou advice
# Rebuild the LR-pairs data table to reflect the full 3x3 cluster combination
# (ECs, FAPs, sCs) instead of the original 2x2 (FAPs, sCs) combination.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$clusters = @("ECs", "FAPs", "sCs")

# Per-sending-cluster ligand stats: F,G,H,I,J
$ligandStats = @{
    "ECs"  = @(0.3333333333333333, 0.071327, 0.213981, 0.03356605248408491, 0.03356605248408491)
    "FAPs" = @(1, 0.308223, 0.9246690000000001, 0.1450478695977975, 0.1450478695977975)
    "sCs"  = @(1, 1.745424333333333, 5.236273, 0.8213860779181176, 0.8213860779181176)
}

# Per-target-cluster receptor stats: K,L,M,N,O,P
$receptorStats = @{
    "ECs"  = @(2, 0.6666666666666666, 2.856403666666667, 8.569210999999999, 0.235832554697756, 0.235832554697756)
    "FAPs" = @(3, 1, 7.000300666666668, 21.000902, 0.5779641054021444, 0.5779641054021444)
    "sCs"  = @(3, 1, 2.255294666666666, 6.765884, 0.1862033399000996, 0.1862033399000996)
}

# Edge weights Q,R,S,T keyed by "sending|target"
$edgeStats = @{
    "ECs|ECs"   = @(0.2037387043323333, 1.833648338991, 0.007915967908440704, 0.007915967908440702)
    "ECs|FAPs"  = @(0.4993104456513334, 4.493794010862, 0.01939997349584556, 0.01939997349584556)
    "ECs|sCs"   = @(0.1608634026893333, 1.447770624204, 0.006250111079798644, 0.006250111079798644)
    "FAPs|ECs"  = @(0.8804093073510001, 7.923683766159, 0.03420700964071556, 0.03420700964071556)
    "FAPs|FAPs" = @(2.157653672382001, 19.41888305143801, 0.08383246219257792, 0.08383246219257792)
    "FAPs|sCs"  = @(0.695133688044, 6.256203192396001, 0.02700839776450401, 0.02700839776450401)
    "sCs|ECs"   = @(4.985636465622555, 44.870728190603, 0.1937095771485997, 0.1937095771485997)
    "sCs|FAPs"  = @(12.21849512424956, 109.966456118246, 0.4747316697137209, 0.4747316697137209)
    "sCs|sCs"   = @(3.936446190036888, 35.428015710332, 0.1529448310557969, 0.1529448310557969)
}

$eStats = @{
    "ECs"  = 1
    "FAPs" = 3
    "sCs"  = 3
}

$row = 2
foreach ($sending in $clusters) {
    foreach ($target in $clusters) {
        $ws.Cells.Item($row, 1).Value = $sending
        $ws.Cells.Item($row, 2).Value = "Efnb3"
        $ws.Cells.Item($row, 3).Value = "Epha4"
        $ws.Cells.Item($row, 4).Value = $target

        $ws.Cells.Item($row, 5).Value = $eStats[$sending]

        $lig = $ligandStats[$sending]
        $ws.Cells.Item($row, 6).Value = $lig[0]
        $ws.Cells.Item($row, 7).Value = $lig[1]
        $ws.Cells.Item($row, 8).Value = $lig[2]
        $ws.Cells.Item($row, 9).Value = $lig[3]
        $ws.Cells.Item($row, 10).Value = $lig[4]

        $rec = $receptorStats[$target]
        $ws.Cells.Item($row, 11).Value = $rec[0]
        $ws.Cells.Item($row, 12).Value = $rec[1]
        $ws.Cells.Item($row, 13).Value = $rec[2]
        $ws.Cells.Item($row, 14).Value = $rec[3]
        $ws.Cells.Item($row, 15).Value = $rec[4]
        $ws.Cells.Item($row, 16).Value = $rec[5]

        $edge = $edgeStats["$sending|$target"]
        $ws.Cells.Item($row, 17).Value = $edge[0]
        $ws.Cells.Item($row, 18).Value = $edge[1]
        $ws.Cells.Item($row, 19).Value = $edge[2]
        $ws.Cells.Item($row, 20).Value = $edge[3]

        $row = $row + 1
    }
}
